$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('ECs','Il1a','Il1rap','ECs',1,0.3333333333333333,94.03756700000001,282.112701,0.948690694112009,0.948690694112009,3,1,41.89758066666666,125.692742,0.7683440336868811,0.768344033686881,3939.946549079571,35459.51894171615,0.728920834635228,0.7289208346352279),
  @('ECs','Il1a','Il1rap','FAPs',1,0.3333333333333333,94.03756700000001,282.112701,0.948690694112009,0.948690694112009,3,1,3.697151666666667,11.091455,0.06780067916854361,0.0678006791685436,347.6711475633284,3129.040328069955,0.06432187338167127,0.06432187338167125),
  @('ECs','Il1a','Il1rap','M2',1,0.3333333333333333,94.03756700000001,282.112701,0.948690694112009,0.948690694112009,3,1,4.857959,14.573877,0.08908828992398354,0.08908828992398353,456.8306449457531,4111.475804511778,0.08451723160523585,0.08451723160523583),
  @('ECs','Il1a','Il1rap','sCs',1,0.3333333333333333,94.03756700000001,282.112701,0.948690694112009,0.948690694112009,3,1,4.077023,12.231069,0.07476699722059184,0.07476699722059182,383.393323523041,3450.539911707369,0.07093075448987392,0.0709307544898739),
  @('M2','Il1a','Il1rap','ECs',3,1,5.085959333333333,15.257878,0.05130930588799102,0.05130930588799102,3,1,41.89758066666666,125.692742,0.7683440336868811,0.768344033686881,213.0893914357195,1917.804522921476,0.03942319905165306,0.03942319905165306),
  @('M2','Il1a','Il1rap','FAPs',3,1,5.085959333333333,15.257878,0.05130930588799102,0.05130930588799102,3,1,3.697151666666667,11.091455,0.06780067916854361,0.0678006791685436,18.80356302583222,169.23206723249,0.003478805786872345,0.003478805786872344),
  @('M2','Il1a','Il1rap','M2',3,1,5.085959333333333,15.257878,0.05130930588799102,0.05130930588799102,3,1,4.857959,14.573877,0.08908828992398354,0.08908828992398353,24.70738191700067,222.366437253006,0.0045710583187477,0.004571058318747699),
  @('M2','Il1a','Il1rap','sCs',3,1,5.085959333333333,15.257878,0.05130930588799102,0.05130930588799102,3,1,4.077023,12.231069,0.07476699722059184,0.07476699722059182,20.73557317906467,186.620158611582,0.003836242730717921,0.00383624273071792)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $rowData = $data[$i]
  for ($j = 0; $j -lt $rowData.Count; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $rowData[$j]
  }
}
